$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.17 = 12165.36 pesos`n✅ 12165.36 pesos = 3.18 = 974.18 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2.Range("N10").Value = 315.5
$ws2.Range("O10").Value = 3838.17
$ws2.Range("N12").Value = 3831.5
$ws2.Range("O12").Value = 306.82
